$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 335, shifting existing rows 335:358 down to 336:359
$ws.Rows.Item(335).EntireRow.Insert()

# Populate the newly inserted row 335 with the new price record
$ws.Range("A335").Value = 10
$ws.Range("B335").Value = "Vega Modelo de Temuco"
$ws.Range("C335").Value = "La Araucanía"
$ws.Range("D335").Value = 44585
$ws.Range("E335").Value = 9
$ws.Range("F335").Value = 100112023
$ws.Range("G335").Value = "Brócoli"
$ws.Range("H335").Value = "Sin especificar"
$ws.Range("I335").Value = "Primera"
$ws.Range("J335").Value = 650
$ws.Range("K335").Value = 1000
$ws.Range("L335").Value = 1000
$ws.Range("M335").Value = 1000
$ws.Range("N335").Value = "$/unidad"
$ws.Range("O335").Value = "Provincia de Cautín"
$ws.Range("P335").Value = 1000
$ws.Range("Q335").Value = 1
$ws.Range("R335").Value = "Hortaliza"
